# Update the GUID-based file names and timestamps as part of regenerating
# the localization handoff/handback report.

$wb = $excel.ActiveWorkbook

$oldGuid = "8d08f052-a844-434b-8ac3-2ddd7c558da3"
$newGuid = "17d42d58-5e2d-4784-9aef-f69d7b3e93c1"

$oldZhCn = "8d08f052-a844-434b-8ac3-2ddd7c558da3.280d71ad9f6a803635c1bba13ce277061b4978e4.zh-cn.xlf"
$newZhCn = "17d42d58-5e2d-4784-9aef-f69d7b3e93c1.85ece963d56a4f108e92a8f30855354f066ec6eb.zh-cn.xlf"

$oldDeDe = "8d08f052-a844-434b-8ac3-2ddd7c558da3.280d71ad9f6a803635c1bba13ce277061b4978e4.de-de.xlf"
$newDeDe = "17d42d58-5e2d-4784-9aef-f69d7b3e93c1.85ece963d56a4f108e92a8f30855354f066ec6eb.de-de.xlf"

# --- Overview sheet ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("A2").Value = "$newGuid.md"
$wsOverview.Range("B2").Value = "e2e\$newGuid.md"
$wsOverview.Hyperlinks.Item(1).TextToDisplay = "e2e\$newGuid.md"
$wsOverview.Range("G2").Value = "2016-09-07 07:21:59"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("A2").Value = "$newGuid.md"
$wsZhCn.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsZhCn.Range("G2").Value = $newZhCn
$wsZhCn.Range("H2").Value = "2016-09-07 07:21:52"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("A2").Value = "$newGuid.md"
$wsDeDe.Hyperlinks.Item(1).TextToDisplay = "$newGuid.md"
$wsDeDe.Range("G2").Value = $newDeDe
$wsDeDe.Range("H2").Value = "2016-09-07 07:21:59"
